$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D4 "Car endpoints:" block ---
# Remove the trailing blank line after "DELETE" and re-apply the
# "headline" run (bold, bigger, green) that matches the "Entity: Car"
# style used elsewhere in the sheet.
$d4 = $ws.Range("D4")
$d4.Value = "Car endpoints:`nPOST`nGET`nGET ALL`nGET BY NAME`nPUT`nDELETE"

$d4Head = $d4.Characters(1, 14)
$d4Head.Font.Bold = $true
$d4Head.Font.Size = 14
$d4Head.Font.Color = 5287936

$d4Len = $d4.Characters().Text.Length
$d4Body = $d4.Characters(15, $d4Len - 14)
$d4Body.Font.Bold = $false
$d4Body.Font.Size = 11

# --- D5 "Documentation:" block ---
# Same headline re-formatting; text itself is unchanged.
$d5 = $ws.Range("D5")
$d5Head = $d5.Characters(1, 14)
$d5Head.Font.Bold = $true
$d5Head.Font.Size = 14
$d5Head.Font.Color = 5287936

$d5Len = $d5.Characters().Text.Length
$d5Body = $d5.Characters(15, $d5Len - 14)
$d5Body.Font.Bold = $false
$d5Body.Font.Size = 11

# --- Row heights grew to fit the larger headline text ---
$ws.Rows.Item(2).RowHeight = 44.25
$ws.Rows.Item(3).RowHeight = 142.5
$ws.Rows.Item(4).RowHeight = 126.75
$ws.Rows.Item(5).RowHeight = 84

# --- Selection moved to E2 ---
$ws.Range("E2").Select() | Out-Null
